$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new "environ" indicator column: set value 1 for E2:E7 (rows under the
# existing "environ" header in E1)
$ws.Range("E2:E7").Value = 1

# Move the active cell selection to E1, matching the updated file state
$ws.Range("E1").Select()
